# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff
# (Thu Aug 24 14:59:42 UTC 2023 GitHub Actions update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, $Ref, $Text)
    $cell = $Sheet.Range($Ref)
    # Force text storage so numeric-looking strings (e.g. "219.54")
    # are not reinterpreted as numbers, matching the source data which
    # stores every value (price/volume/coin/link) as plain text.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" "26.276.23"
Set-CellText $ws "E2" "  +0.71%  "
Set-CellText $ws "D3" "1.657.20"
Set-CellText $ws "E3" "  +0.10%  "
Set-CellText $ws "E4" "  +0.51%  "
Set-CellText $ws "D5" "219.54"
Set-CellText $ws "E5" "  +2.21%  "
Set-CellText $ws "D6" "0.5217"
Set-CellText $ws "E6" "  -0.55%  "
Set-CellText $ws "D7" "1.005"
Set-CellText $ws "E7" "  +0.49%  "
Set-CellText $ws "E8" "  +0.91%  "
Set-CellText $ws "D9" "0.06328"
Set-CellText $ws "E9" "  -0.91%  "
Set-CellText $ws "D10" "21.37"
Set-CellText $ws "E10" "  +2.63%  "
Set-CellText $ws "D11" "0.07761"
Set-CellText $ws "E11" "  +0.26%  "
Set-CellText $ws "B12" "Polkadot"
Set-CellText $ws "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText $ws "D12" "4.443"
Set-CellText $ws "E12" "  -0.37%  "
Set-CellText $ws "B13" "WrappedEther"
Set-CellText $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-CellText $ws "D13" "1.639.74"
Set-CellText $ws "E13" "  -0.89%  "
Set-CellText $ws "D14" "0.5480"
Set-CellText $ws "E14" "  -0.65%  "
Set-CellText $ws "D15" "0.0₅8235"
Set-CellText $ws "E15" "  -0.86%  "
Set-CellText $ws "D16" "64.98"
Set-CellText $ws "E16" "  -0.16%  "
Set-CellText $ws "D17" "26.276.58"
Set-CellText $ws "E17" "  +0.62%  "
Set-CellText $ws "D18" "1.006"
Set-CellText $ws "E18" "  +0.47%  "
Set-CellText $ws "D19" "4.696"
Set-CellText $ws "E19" "  -1.25%  "
Set-CellText $ws "D20" "191.65"
Set-CellText $ws "E20" "  +0.62%  "
Set-CellText $ws "D21" "10.20"
Set-CellText $ws "E21" "  -0.58%  "
Set-CellText $ws "D22" "6.216"
Set-CellText $ws "E22" "  -2.20%  "
Set-CellText $ws "E23" "  +0.65%  "
Set-CellText $ws "D24" "138.99"
Set-CellText $ws "E24" "  -2.85%  "
Set-CellText $ws "D25" "0.1255"
Set-CellText $ws "E25" "  +0.34%  "
Set-CellText $ws "D26" "7.305"
Set-CellText $ws "E26" "  -1.49%  "
Set-CellText $ws "D27" "16.09"
Set-CellText $ws "E27" "  +0.30%  "
Set-CellText $ws "D28" "1.420"
Set-CellText $ws "E28" "  +0.06%  "
Set-CellText $ws "D29" "0.06057"
Set-CellText $ws "E29" "  +1.79%  "
Set-CellText $ws "E30" "  +2.22%  "
Set-CellText $ws "D31" "3.556"
Set-CellText $ws "E31" "  +3.36%  "
Set-CellText $ws "D32" "3.375"
Set-CellText $ws "E32" "  -1.25%  "
Set-CellText $ws "D33" "1.663"
Set-CellText $ws "E33" "  +0.42%  "
Set-CellText $ws "D34" "0.9888"
Set-CellText $ws "E34" "  -0.98%  "
Set-CellText $ws "D35" "2.427"
Set-CellText $ws "E35" "  +1.19%  "
Set-CellText $ws "D36" "2.773"
Set-CellText $ws "E36" "  +0.29%  "
Set-CellText $ws "E37" "  +5.99%  "
Set-CellText $ws "D38" "0.01597"
Set-CellText $ws "E38" "  -0.52%  "
Set-CellText $ws "D39" "5.983"
Set-CellText $ws "E39" "  +1.80%  "
Set-CellText $ws "D40" "1.070.87"
Set-CellText $ws "E40" "  +4.21%  "
Set-CellText $ws "D41" "0.8514"
Set-CellText $ws "E41" "  -0.65%  "
Set-CellText $ws "E42" "  +0.42%  "
Set-CellText $ws "D43" "99.98"
Set-CellText $ws "E43" "  +0.55%  "
Set-CellText $ws "D44" "1.799.99"
Set-CellText $ws "E44" "  -0.14%  "
Set-CellText $ws "B45" "BabyDogeCoin"
Set-CellText $ws "C45" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText $ws "D45" "0.0₈110"
Set-CellText $ws "E45" "  +2.45%  "
Set-CellText $ws "B46" "Aave"
Set-CellText $ws "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText $ws "D46" "57.55"
Set-CellText $ws "E46" "  +2.87%  "
Set-CellText $ws "E47" "  +0.16%  "
Set-CellText $ws "D48" "8.052"
Set-CellText $ws "E48" "  -0.23%  "
Set-CellText $ws "D49" "0.05192"
Set-CellText $ws "E49" "  +0.71%  "
Set-CellText $ws "D50" "1.474"
Set-CellText $ws "E50" "  +5.80%  "
Set-CellText $ws "D51" "0.4232"
Set-CellText $ws "E51" "  +0.56%  "
